$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 96
$ws.Range("F5").Value = 407
$ws.Range("F6").Value = 1810
$ws.Range("F7").Value = 854
$ws.Range("F8").Value = 1265
$ws.Range("F9").Value = 66
$ws.Range("F10").Value = 436
$ws.Range("F11").Value = 127
$ws.Range("F12").Value = 2619
$ws.Range("F13").Value = 351
$ws.Range("F14").Value = 862
$ws.Range("F15").Value = 1067
$ws.Range("F17").Value = 14
$ws.Range("F19").Value = 1542
$ws.Range("F20").Value = 15
$ws.Range("F21").Value = 1219
$ws.Range("F22").Value = 162
$ws.Range("F24").Value = 1367
$ws.Range("C25").Value = "上海·幻想乡动漫游戏节"
$ws.Range("D25").Value = "中山北路3300号 上海JOYPOLIS世嘉都市乐园"
$ws.Range("E25").Value = "2024.07.13 10:00-07.14 17:00"
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=87440"
$ws.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202406/kJiNcksB1719222316678.jpeg"
$ws.Range("C26").Value = "上海·恋与深空only【女生专场】"
$ws.Range("D26").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws.Range("E26").Value = "2024.07.13 11:00-07.13 21:00"
$ws.Range("F26").Value = 1384
$ws.Range("G26").Value = 88
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=87346"
$ws.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202406/vaE8H0CC1718081128645.jpeg"
$ws.Range("C27").Value = "上海·第五届燃梦BACG PRO动漫嘉年华·我们在燃梦相遇吧！"
$ws.Range("D27").Value = "盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)"
$ws.Range("E27").Value = "2024.07.13 11:00-07.14 16:00"
$ws.Range("F27").Value = 944
$ws.Range("G27").Value = 65.8
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=85235"
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202405/A9VkjJzh1715108077210.jpeg"
$ws.Range("F28").Value = 1301
$ws.Range("F29").Value = 193
$ws.Range("F30").Value = 1256
$ws.Range("F31").Value = 417
$ws.Range("F35").Value = 1814
$ws.Range("F36").Value = 453
$ws.Range("F37").Value = 32
$ws.Range("F39").Value = 17
$ws.Range("F40").Value = 2232
$ws.Range("F41").Value = 127
$ws.Range("F42").Value = 876
$ws.Range("F43").Value = 2732
$ws.Range("F44").Value = 9
$ws.Range("F46").Value = 142

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 56
$ws.Range("F7").Value = 26
$ws.Range("F11").Value = 14
$ws.Range("F17").Value = 58
$ws.Range("F18").Value = 58
$ws.Range("F20").Value = 282
$ws.Range("F26").Value = 57
$ws.Range("F27").Value = 57
$ws.Range("F29").Value = 42
$ws.Range("F30").Value = 210
$ws.Range("F34").Value = 83
$ws.Range("F36").Value = 160

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 3005
$ws.Range("F6").Value = 4834
$ws.Range("F7").Value = 173
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 648
$ws.Range("F10").Value = 907
$ws.Range("F11").Value = 530
$ws.Range("F12").Value = 580
$ws.Range("F13").Value = 1319
$ws.Range("F14").Value = 372
$ws.Range("F15").Value = 1133

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 4834
$ws.Range("F6").Value = 648
$ws.Range("F7").Value = 907
$ws.Range("F8").Value = 530
$ws.Range("F9").Value = 580
$ws.Range("F10").Value = 1319
$ws.Range("F11").Value = 407
$ws.Range("F12").Value = 1811
$ws.Range("F13").Value = 854
$ws.Range("F14").Value = 1265
$ws.Range("F15").Value = 26
$ws.Range("F16").Value = 436
$ws.Range("F17").Value = 1133
$ws.Range("F18").Value = 2619
$ws.Range("F20").Value = 351
$ws.Range("F21").Value = 862
$ws.Range("F22").Value = 1067
$ws.Range("F24").Value = 1542
$ws.Range("F25").Value = 14
$ws.Range("F27").Value = 1219
$ws.Range("F28").Value = 162
$ws.Range("F30").Value = 1385
$ws.Range("F31").Value = 944
$ws.Range("F32").Value = 1301
$ws.Range("F33").Value = 193
$ws.Range("F35").Value = 58
$ws.Range("F36").Value = 1256
$ws.Range("F37").Value = 417
$ws.Range("F40").Value = 1814
$ws.Range("F41").Value = 57
$ws.Range("F42").Value = 32
$ws.Range("F44").Value = 2232
$ws.Range("F45").Value = 127
$ws.Range("F46").Value = 876
$ws.Range("F47").Value = 2732
$ws.Range("F49").Value = 142
